$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formula")

# Delete the entire row 33 (FC_NET_INCOME_MINUS_CASH_DIV_TOTAL_EQUITY_BNK),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(33).Delete()

# Reflect the resulting selection state (row 33, which now holds the
# content that used to be row 34).
$ws.Range("A33:XFD33").Select()
